$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row text: "N вит." -> "N" and "R оц." -> "R"
$ws.Range("B1").Value = "N"
$ws.Range("C1").Value = "R"

# Update selected cell from K12 to E4
$ws.Range("E4").Select()
